$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'hindern',
    'liegen',
    'lockern',
    'fischen',
    'räumen',
    'pflegen',
    'stehlen',
    'spielen',
    'ruhen',
    'heben',
    'schwächen',
    'sagen',
    'fordern',
    'biegen',
    'saufen',
    'wiegen',
    'wüten',
    'zeugen',
    'lesen',
    'dienen',
    'gründen',
    'knien',
    'tropfen',
    'tanzen',
    'jagen',
    'töten',
    'hauen',
    'klettern',
    'narren',
    'bitten',
    'führen',
    'reizen',
    'trauen',
    'kranken',
    'schwören',
    'kriegen',
    'posten',
    'betteln',
    'wenden',
    'backen',
    'planen',
    'fließen',
    'reiten',
    'zählen',
    'freuen',
    'suchen',
    'arten',
    'bellen',
    'irren',
    'mühen',
    'ehren',
    'füttern',
    'knarren',
    'wandern',
    'wellen',
    'pflanzen',
    'kümmern',
    'tauschen',
    'wehtun',
    'fallen',
    'kehren',
    'warnen',
    'siegen',
    'meistern',
    'wundern',
    'bauen',
    'werfen',
    'sorgen',
    'fesseln',
    'süßen',
    'helfen',
    'mauern',
    'grüßen',
    'schmecken',
    'werden',
    'bluten',
    'schätzen',
    'fassen',
    'sinken',
    'schrecken',
    'formen',
    'flüchten',
    'tollen',
    'krachen',
    'fällen',
    'lohnen',
    'feiern',
    'liefern',
    'greifen',
    'lächeln',
    'hören',
    'malen',
    'fehlen',
    'trennen',
    'quälen',
    'wehren',
    'wärmen',
    'streichen',
    'decken',
    'sterben',
    'flehen',
    'pressen',
    'fahren',
    'fangen',
    'achten',
    'klingen',
    'zahlen',
    'heilen',
    'schlafen',
    'scheinen',
    'erben',
    'wehen',
    'brauchen',
    'boxen',
    'zünden',
    'prügeln',
    'sperren',
    'schleppen',
    'drucken',
    'stecken',
    'starren',
    'stören',
    'gelten',
    'kürzen',
    'stoppen',
    'mögen',
    'loben',
    'schützen',
    'schulden',
    'äußern',
    'kosten',
    'pfeifen',
    'rasen',
    'orten',
    'zeigen',
    'schlucken',
    'stammen',
    'lenken',
    'folgen',
    'kichern',
    'dürfen',
    'wachsen',
    'stechen',
    'spüren',
    'platzen',
    'runden',
    'ändern',
    'fühlen',
    'ärgern',
    'treiben',
    'drehen',
    'enden',
    'schreiten',
    'schlagen',
    'herzen',
    'proben',
    'altern',
    'haben',
    'öffnen',
    'zielen',
    'rauchen',
    'regnen',
    'dringen',
    'sprengen',
    'wirken',
    'schenken',
    'wohnen',
    'sichern',
    'schalten',
    'filmen',
    'treten',
    'weichen',
    'zögern',
    'scheitern',
    'geben',
    'machen',
    'nutzen',
    'rufen',
    'heulen',
    'fügen',
    'schwingen',
    'bremsen',
    'spinnen',
    'bergen',
    'kämpfen',
    'beten',
    'seufzen',
    'sprechen',
    'münzen',
    'graben',
    'jubeln',
    'lügen'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}
